# Update the task list on Hárok1 to reflect final project status:
#  - Tomáš Adam's tasks get expanded scope (UI & deployment, preprocessing, data/task list)
#  - Lucia Szalonová gains an "App Deployment / Major debugging & Final cleanup" task (row 16)
#  - Kristián Maťašovský's "Participants" task moves down to the now-empty row 17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Tomáš Adam / User Interface -> User Interface & App deployment
$ws.Range("D3").Value = "User Interface & App deployment"
$ws.Range("E3").Value = "Project architecture, Repository maintaining, Major Final debuging & cleanup, App finalization & deployment"

# Row 4 - Tomáš Adam / Preprocessing -> Data Preprocessing
$ws.Range("D4").Value = "Data Preprocessing"
$ws.Range("E4").Value = "Column names correction, Value trimming \ factorization & translation, Missing values replacement "

# Row 5 - Tomáš Adam / UI - Data -> UI - Data & Task List
$ws.Range("D5").Value = "UI - Data & Task List"
$ws.Range("E5").Value = "Data tables & Tasks visualizaion"

# Row 15 - Lucia Szalonová / Bid progress -> UI - Bid progress & Data Preprocessing
$ws.Range("D15").Value = "UI - Bid progress & Data Preprocessing"

# Row 16 - was Kristián Maťašovský / Participants, now Lucia Szalonová / App Deployment
$ws.Range("C16").Value = "Lucia Szalonová"
$ws.Range("D16").Value = "App Deployment"
$ws.Range("E16").Value = "Major debuging & Final cleanup"

# Row 17 - now holds what used to be row 16 (Kristián Maťašovský / Participants)
$ws.Range("C17").Value = "Kristián Maťašovský"
$ws.Range("D17").Value = "Participants"
$ws.Range("E17").Value = "Visualization of most succesful participants according to type and category of auction."

# Keep current selection on E8 as in the final saved state
$ws.Range("E8").Select()
